$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: new shuffled permutation of 0..19 (rows 1-20)
$bValues = @(6, 11, 15, 19, 10, 18, 8, 16, 13, 4, 1, 17, 14, 5, 2, 9, 3, 7, 12, 0)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# D1, D2: updated fitness/penalty improvement values
$ws.Range("D1").Value = 112.9106100093386
$ws.Range("D2").Value = 68.86544096462683

# B21: updated last generation fit value
$ws.Range("B21").Value = 0.7981402113560117
